$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: "List 1" summary header block (A36:F36) ---
# Write the text cells in the same order the source workbook's edit used,
# so new shared-string entries land at the same indices as the target file
# (Cue - Med, List 1, Med -Med, Med - Target, List 2).
$ws.Range("C36").Value = "Cue - Med"
$ws.Range("A36").Value = "List 1"
$ws.Range("D36").Value = "Med -Med"
$ws.Range("E36").Value = "Med - Target"
$ws.Range("A36:F36").Font.Bold = $true

# --- Rows 37-40: Mean / SD / Min / Max for List 1 (rows 3:17) ---
$ws.Range("B37").Value = "Mean"
$ws.Range("B37").Font.Bold = $true
$ws.Range("C37").Formula = "=AVERAGE(E3:E17)"
$ws.Range("D37").Formula = "=AVERAGE(I3:I17)"
$ws.Range("E37").Formula = "=AVERAGE(M3:M17)"

$ws.Range("B38").Value = "SD"
$ws.Range("B38").Font.Bold = $true
$ws.Range("C38").Formula = "=STDEV(E3:E17)"
$ws.Range("D38").Formula = "=STDEV(F18:I32)"
$ws.Range("E38").Formula = "=STDEV(M3:M17)"

$ws.Range("B39").Value = "Min"
$ws.Range("B39").Font.Bold = $true
$ws.Range("C39").Formula = "=MIN(E3:E17)"
$ws.Range("D39").Formula = "=MIN(I3:I17)"
$ws.Range("E39").Formula = "=MIN(M3:M17)"

$ws.Range("B40").Value = "Max"
$ws.Range("B40").Font.Bold = $true
$ws.Range("C40").Formula = "=MAX(E3:E17)"
$ws.Range("D40").Formula = "=MAX(I3:I17)"
$ws.Range("E40").Formula = "=MAX(M3:M17)"

# --- Row 42: "List 2" summary header block (A42:E42) ---
$ws.Range("C42").Value = "Cue - Med"
$ws.Range("D42").Value = "Med -Med"
$ws.Range("E42").Value = "Med - Target"
$ws.Range("A42").Value = "List 2"
$ws.Range("A42:E42").Font.Bold = $true

# --- Rows 43-46: Mean / SD / Min / Max for List 2 (rows 18:32) ---
$ws.Range("B43").Value = "Mean"
$ws.Range("B43").Font.Bold = $true
$ws.Range("C43").Formula = "=AVERAGE(E18:E32)"
$ws.Range("D43").Formula = "=AVERAGE(I18:I32)"
$ws.Range("E43").Formula = "=AVERAGE(M18:M32)"

$ws.Range("B44").Value = "SD"
$ws.Range("B44").Font.Bold = $true
$ws.Range("C44").Formula = "=STDEV(E18:E32)"
$ws.Range("D44").Formula = "=STDEV(F18:I32)"
$ws.Range("E44").Formula = "=STDEV(M18:M32)"

$ws.Range("B45").Value = "Min"
$ws.Range("B45").Font.Bold = $true
$ws.Range("C45").Formula = "=MIN(E18:E32)"
$ws.Range("D45").Formula = "=MIN(I18:I32)"
$ws.Range("E45").Formula = "=MIN(M18:M32)"

$ws.Range("B46").Value = "Max"
$ws.Range("B46").Font.Bold = $true
$ws.Range("C46").Formula = "=MAX(E18:E32)"
$ws.Range("D46").Formula = "=MAX(I18:I32)"
$ws.Range("E46").Formula = "=MAX(M18:M32)"

# --- View state: scroll + selection matches the saved file ---
[void]$ws.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 22
[void]$ws.Range("E47").Select()
